# Scheduled runner update: refresh computed profit columns (H-N) on the
# per-class "leve" sheets to reflect the latest market data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 17340.5
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = $null

$ws.Range("H23").Value = 17340.5
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null

$ws.Range("H57").Value = 37483.75
$ws.Range("I57").Value = 25000
$ws.Range("J57").Value = 49967.5
$ws.Range("K57").Value = 75000
$ws.Range("L57").Value = 149902.5
$ws.Range("M57").Value = -74501
$ws.Range("N57").Value = -150900.5

$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = $null

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = $null

$ws.Range("H93").Value = 60783.668
$ws.Range("J93").Value = 60783.668
$ws.Range("L93").Value = 60783.668
$ws.Range("N93").Value = -65775.66800000001

$ws.Range("H132").Value = 37037.55
$ws.Range("I132").Value = 52795.05
$ws.Range("J132").Value = 2020.8889
$ws.Range("K132").Value = 158385.15
$ws.Range("L132").Value = 6062.6667
$ws.Range("M132").Value = -155855.15
$ws.Range("N132").Value = -11122.6667

$ws.Range("H135").Value = 590.8333
$ws.Range("I135").Value = 652.2857
$ws.Range("K135").Value = 5870.571300000001
$ws.Range("M135").Value = -3335.571300000001

$ws.Range("H141").Value = 4433.6665
$ws.Range("I141").Value = 3621.2
$ws.Range("K141").Value = 10863.6
$ws.Range("M141").Value = -5683.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 50110
$ws.Range("J56").Value = 50110
$ws.Range("L56").Value = 50110
$ws.Range("N56").Value = -51594

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 4152.8335
$ws.Range("I5").Value = 3983.6
$ws.Range("K5").Value = 3983.6
$ws.Range("M5").Value = -3870.6

$ws.Range("H22").Value = 569.25
$ws.Range("I22").Value = 550
$ws.Range("J22").Value = 571
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 571
$ws.Range("M22").Value = -377
$ws.Range("N22").Value = -917

$ws.Range("H94").Value = 1720.3636
$ws.Range("I94").Value = 1720.3636
$ws.Range("K94").Value = 1720.3636
$ws.Range("M94").Value = -1269.3636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2252.5833
$ws.Range("I31").Value = 1285.3572
$ws.Range("K31").Value = 1285.3572
$ws.Range("M31").Value = -990.3571999999999

$ws.Range("H34").Value = 2252.5833
$ws.Range("I34").Value = 1285.3572
$ws.Range("K34").Value = 1285.3572
$ws.Range("M34").Value = -1083.3572

$ws.Range("H58").Value = 3074.4
$ws.Range("I58").Value = 3141.25
$ws.Range("J58").Value = 2974.125
$ws.Range("K58").Value = 3141.25
$ws.Range("L58").Value = 2974.125
$ws.Range("M58").Value = -2938.25
$ws.Range("N58").Value = -3380.125

$ws.Range("H109").Value = 37993.5
$ws.Range("J109").Value = 44988
$ws.Range("L109").Value = 44988
$ws.Range("N109").Value = -47068

$ws.Range("H122").Value = 29988.445
$ws.Range("I122").Value = 2599.3333
$ws.Range("K122").Value = 7797.999899999999
$ws.Range("M122").Value = -5347.999899999999

$ws.Range("H134").Value = 2285.923
$ws.Range("I134").Value = 2221.9
$ws.Range("K134").Value = 6665.700000000001
$ws.Range("M134").Value = -4130.700000000001

$ws.Range("H136").Value = 3074.4
$ws.Range("I136").Value = 3141.25
$ws.Range("J136").Value = 2974.125
$ws.Range("K136").Value = 9423.75
$ws.Range("L136").Value = 8922.375
$ws.Range("M136").Value = -6873.75
$ws.Range("N136").Value = -14022.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 488.2143
$ws.Range("J12").Value = 635.2222
$ws.Range("L12").Value = 1905.6666
$ws.Range("N12").Value = -2251.6666

$ws.Range("H23").Value = 2311.5
$ws.Range("J23").Value = 1865.9412
$ws.Range("L23").Value = 5597.8236
$ws.Range("N23").Value = -6067.8236

$ws.Range("H33").Value = 119.28571
$ws.Range("J33").Value = 149.66667
$ws.Range("L33").Value = 898.0000200000001
$ws.Range("N33").Value = -1464.00002

$ws.Range("H82").Value = 1500
$ws.Range("I82").Value = 1500
$ws.Range("K82").Value = 4500
$ws.Range("M82").Value = -4094

$ws.Range("H85").Value = 1500
$ws.Range("I85").Value = 1500
$ws.Range("K85").Value = 4500
$ws.Range("M85").Value = -3096

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8750.5
$ws.Range("I5").Value = 5004
$ws.Range("J5").Value = 19990
$ws.Range("K5").Value = 5004
$ws.Range("L5").Value = 19990
$ws.Range("M5").Value = -4892
$ws.Range("N5").Value = -20214

$ws.Range("H70").Value = 4647.8887
$ws.Range("I70").Value = 3867.5386
$ws.Range("J70").Value = 6676.8
$ws.Range("K70").Value = 3867.5386
$ws.Range("L70").Value = 6676.8
$ws.Range("M70").Value = -3597.5386
$ws.Range("N70").Value = -7216.8

$ws.Range("H73").Value = 4647.8887
$ws.Range("I73").Value = 3867.5386
$ws.Range("J73").Value = 6676.8
$ws.Range("K73").Value = 3867.5386
$ws.Range("L73").Value = 6676.8
$ws.Range("M73").Value = -2931.5386
$ws.Range("N73").Value = -8548.799999999999

$ws.Range("H116").Value = 79897.5
$ws.Range("J116").Value = 79897.5
$ws.Range("L116").Value = 79897.5
$ws.Range("N116").Value = -89075.5

$ws.Range("H132").Value = 2812
$ws.Range("I132").Value = 2374.4167
$ws.Range("K132").Value = 7123.250100000001
$ws.Range("M132").Value = -4593.250100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 22244
$ws.Range("J11").Value = 22244
$ws.Range("L11").Value = 22244
$ws.Range("N11").Value = -22524

$ws.Range("H22").Value = 3666.111
$ws.Range("J22").Value = 5199.8
$ws.Range("L22").Value = 5199.8
$ws.Range("N22").Value = -5789.8

$ws.Range("H25").Value = 12555.444
$ws.Range("I25").Value = 6665.3335
$ws.Range("J25").Value = 24335.666
$ws.Range("K25").Value = 6665.3335
$ws.Range("L25").Value = 24335.666
$ws.Range("M25").Value = -6435.3335
$ws.Range("N25").Value = -24795.666

$ws.Range("H27").Value = 3666.111
$ws.Range("J27").Value = 5199.8
$ws.Range("L27").Value = 5199.8
$ws.Range("N27").Value = -5413.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 22426.875
$ws.Range("J70").Value = 22426.875
$ws.Range("L70").Value = 22426.875
$ws.Range("N70").Value = -23056.875

$ws.Range("H73").Value = 22426.875
$ws.Range("J73").Value = 22426.875
$ws.Range("L73").Value = 22426.875
$ws.Range("N73").Value = -24610.875

$ws.Range("H109").Value = 20062.5
$ws.Range("J109").Value = 20062.5
$ws.Range("L109").Value = 20062.5
$ws.Range("N109").Value = -22836.5
